# Add team record (Wins/Losses/Ties) columns to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine last used row (data rows 2..51, header row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 51 }

# Headers in row 1 for columns AD, AE, AF - match style of existing header cells (copy from AC1)
$headerStyleSource = $ws.Range("AC1")

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerStyleSource.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill data rows with win/loss/tie record
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 94
    $ws.Cells.Item($r, 31).Value = 68
    $ws.Cells.Item($r, 32).Value = 0
}

$wb.Save()
